$wb = $excel.ActiveWorkbook

# "Metadata" worksheet holds the ValueSet-level properties
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-opt-out-reason"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# "Include from Engagement Opt-O" worksheet holds the CodeSystem properties
$codes = $wb.Worksheets.Item("Include from Engagement Opt-O")
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/eng-opt-out-reason"
